$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.08410299302171197
$ws.Range("C2").Value = 1.253677763513526
$ws.Range("D2").Value = 6.261918074916323
$ws.Range("E2").Value = 2.502382479741321
$ws.Range("F2").Value = 2.531284644362528
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = -0.01946844920150798
$ws.Range("C3").Value = 1.273117676343247
$ws.Range("D3").Value = 6.739870713701541
$ws.Range("E3").Value = 2.596126097419296
$ws.Range("F3").Value = 2.628303442435235
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.0393246434631094
$ws.Range("C4").Value = 1.235849058089852
$ws.Range("D4").Value = 6.613468431034368
$ws.Range("E4").Value = 2.571666469632944
$ws.Range("F4").Value = 2.604123358268155
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = -0.04878613897638684
$ws.Range("C5").Value = 1.266707606888937
$ws.Range("D5").Value = 6.697608253045909
$ws.Range("E5").Value = 2.587973773639507
$ws.Range("F5").Value = 2.621339043423467
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.01069592955521734
$ws.Range("C6").Value = 1.306387404899978
$ws.Range("D6").Value = 6.984373269163457
$ws.Range("E6").Value = 2.642796486520189
$ws.Range("F6").Value = 2.678249917277403
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = -0.03495980545002601
$ws.Range("C7").Value = 1.341848499270465
$ws.Range("D7").Value = 6.976450903423854
$ws.Range("E7").Value = 2.64129720088896
$ws.Range("F7").Value = 2.677496044559739
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.04415342685500446
$ws.Range("C8").Value = 1.313189931210337
$ws.Range("D8").Value = 7.030782837840379
$ws.Range("E8").Value = 2.651562339044733
$ws.Range("F8").Value = 2.688802172282798
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = -0.01865779601231733
$ws.Range("C9").Value = 1.387484845158188
$ws.Range("D9").Value = 7.334554756911054
$ws.Range("E9").Value = 2.708238312429513
$ws.Range("F9").Value = 2.747711521186984
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.0924138300119391
$ws.Range("C10").Value = 1.344817476870875
$ws.Range("D10").Value = 7.384053110940232
$ws.Range("E10").Value = 2.717361424422639
$ws.Range("F10").Value = 2.756630766477241
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = -0.01024319886757436
$ws.Range("C11").Value = 1.413099074990776
$ws.Range("D11").Value = 7.676407440777855
$ws.Range("E11").Value = 2.770633039718153
$ws.Range("F11").Value = 2.813571923966318
$ws.Range("G11").Value = 33

